# Update the cryptos price/volume table (columns D and E) for rows 2-51
# to reflect the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text, preventing Excel from
# auto-converting numeric-looking strings (e.g. "1.00" -> 1, "5.80" -> 5.8)
# while leaving the cells original style untouched.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "59.149.61"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.519.20"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.31%  "
Set-TextValue $ws.Range("D5") "535.72"
$ws.Range("E5").Value = "  -0.16%  "
Set-TextValue $ws.Range("D6") "140.19"
$ws.Range("E6").Value = "  -3.41%  "
$ws.Range("E7").Value = "  +0.33%  "
Set-TextValue $ws.Range("D8") "0.563"
$ws.Range("E8").Value = "  -2.07%  "
$ws.Range("D9").Value = "2.525.60"
$ws.Range("E9").Value = "  +0.33%  "
Set-TextValue $ws.Range("D10") "0.0990"
$ws.Range("E10").Value = "  -0.65%  "
Set-TextValue $ws.Range("D11") "0.160"
$ws.Range("E11").Value = "  +1.70%  "
Set-TextValue $ws.Range("D12") "5.39"
$ws.Range("E12").Value = "  -2.69%  "
Set-TextValue $ws.Range("D13") "0.354"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "2.967.73"
$ws.Range("E14").Value = "  +1.40%  "
Set-TextValue $ws.Range("D15") "23.14"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("D16").Value = "59.130.40"
$ws.Range("E16").Value = "  +0.32%  "
Set-TextValue $ws.Range("D17") "0.0000140"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "2.518.69"
$ws.Range("E18").Value = "  +0.13%  "
Set-TextValue $ws.Range("D19") "10.96"
$ws.Range("E19").Value = "  -2.88%  "
Set-TextValue $ws.Range("D20") "4.22"
$ws.Range("E20").Value = "  -1.69%  "
Set-TextValue $ws.Range("D21") "320.23"
$ws.Range("E21").Value = "  -1.10%  "
Set-TextValue $ws.Range("D22") "0.998"
$ws.Range("E22").Value = "  -0.10%  "
Set-TextValue $ws.Range("D23") "5.80"
$ws.Range("E23").Value = "  +0.81%  "
Set-TextValue $ws.Range("D24") "62.44"
$ws.Range("E24").Value = "  +1.67%  "
Set-TextValue $ws.Range("D25") "0.418"
$ws.Range("E25").Value = "  -4.87%  "
$ws.Range("E26").Value = "  +2.71%  "
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.38%  "
Set-TextValue $ws.Range("D28") "7.81"
$ws.Range("E28").Value = "  +0.36%  "
Set-TextValue $ws.Range("D29") "6.76"
$ws.Range("E29").Value = "  -1.49%  "
Set-TextValue $ws.Range("D30") "0.0₃0769"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("E31").Value = "  +0.62%  "
Set-TextValue $ws.Range("D32") "161.91"
$ws.Range("E32").Value = "  +1.59%  "
Set-TextValue $ws.Range("D33") "0.999"
$ws.Range("E33").Value = "  +0.30%  "
Set-TextValue $ws.Range("D34") "1.13"
$ws.Range("E34").Value = "  -9.71%  "
Set-TextValue $ws.Range("D35") "1.44"
$ws.Range("E35").Value = "  -0.36%  "
Set-TextValue $ws.Range("D36") "18.44"
$ws.Range("E36").Value = "  -0.54%  "
Set-TextValue $ws.Range("D37") "4.21"
$ws.Range("E37").Value = "  -5.31%  "
Set-TextValue $ws.Range("D38") "1.58"
$ws.Range("E38").Value = "  -2.13%  "
Set-TextValue $ws.Range("D39") "36.92"
$ws.Range("E39").Value = "  +0.35%  "
Set-TextValue $ws.Range("D40") "3.64"
$ws.Range("E40").Value = "  -1.18%  "
Set-TextValue $ws.Range("D41") "5.32"
$ws.Range("E41").Value = "  -9.84%  "
Set-TextValue $ws.Range("D42") "286.10"
$ws.Range("E42").Value = "  -7.55%  "
Set-TextValue $ws.Range("D43") "0.803"
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("E44").Value = "  +0.06%  "
Set-TextValue $ws.Range("D45") "0.599"
$ws.Range("E45").Value = "  +0.69%  "
Set-TextValue $ws.Range("D46") "10.86"
$ws.Range("E46").Value = "  +0.70%  "
Set-TextValue $ws.Range("D47") "124.61"
$ws.Range("E47").Value = "  +0.08%  "
Set-TextValue $ws.Range("D48") "0.0928"
$ws.Range("E48").Value = "  -0.33%  "
Set-TextValue $ws.Range("D49") "18.60"
$ws.Range("E49").Value = "  -0.01%  "
Set-TextValue $ws.Range("D50") "0.0508"
$ws.Range("E50").Value = "  -1.69%  "
Set-TextValue $ws.Range("D51") "0.0222"
$ws.Range("E51").Value = "  -2.31%  "
